$d = $word.ActiveDocument
$d.Content.Find.Execute("63×54=", $true, $false, $false, $false, $false, $true, 1, $false, "70×28=", 2)
$d.Content.Find.Execute("80×16=", $true, $false, $false, $false, $false, $true, 1, $false, "66×14=", 2)
$d.Content.Find.Execute("34×63=", $true, $false, $false, $false, $false, $true, 1, $false, "59×60=", 2)
$d.Content.Find.Execute("77×19=", $true, $false, $false, $false, $false, $true, 1, $false, "16×14=", 2)
$d.Content.Find.Execute("26×86=", $true, $false, $false, $false, $false, $true, 1, $false, "35×88=", 2)
$d.Content.Find.Execute("66×17=", $true, $false, $false, $false, $false, $true, 1, $false, "65×54=", 2)
$d.Content.Find.Execute("70×88=", $true, $false, $false, $false, $false, $true, 1, $false, "61×23=", 2)
$d.Content.Find.Execute("32×16=", $true, $false, $false, $false, $false, $true, 1, $false, "22×50=", 2)
$d.Content.Find.Execute("27×92=", $true, $false, $false, $false, $false, $true, 1, $false, "84×61=", 2)
$d.Content.Find.Execute("46×34=", $true, $false, $false, $false, $false, $true, 1, $false, "44×60=", 2)
$d.Content.Find.Execute("69×96=", $true, $false, $false, $false, $false, $true, 1, $false, "60×73=", 2)
$d.Content.Find.Execute("90×62=", $true, $false, $false, $false, $false, $true, 1, $false, "26×34=", 2)
$d.Content.Find.Execute("74×90=", $true, $false, $false, $false, $false, $true, 1, $false, "96×12=", 2)
$d.Content.Find.Execute("81×42=", $true, $false, $false, $false, $false, $true, 1, $false, "54×55=", 2)
$d.Content.Find.Execute("63×89=", $true, $false, $false, $false, $false, $true, 1, $false, "54×45=", 2)
$d.Content.Find.Execute("96×63=", $true, $false, $false, $false, $false, $true, 1, $false, "82×81=", 2)
$d.Content.Find.Execute("28×16=", $true, $false, $false, $false, $false, $true, 1, $false, "24×52=", 2)
$d.Content.Find.Execute("99×54=", $true, $false, $false, $false, $false, $true, 1, $false, "20×25=", 2)
$d.Content.Find.Execute("64×61=", $true, $false, $false, $false, $false, $true, 1, $false, "41×89=", 2)
$d.Content.Find.Execute("69×73=", $true, $false, $false, $false, $false, $true, 1, $false, "18×41=", 2)
$d.Content.Find.Execute("78×99=", $true, $false, $false, $false, $false, $true, 1, $false, "38×99=", 2)
$d.Content.Find.Execute("62×38=", $true, $false, $false, $false, $false, $true, 1, $false, "37×94=", 2)
$d.Content.Find.Execute("37×93=", $true, $false, $false, $false, $false, $true, 1, $false, "13×69=", 2)
$d.Content.Find.Execute("53×85=", $true, $false, $false, $false, $false, $true, 1, $false, "68×15=", 2)
$d.Content.Find.Execute("93×56=", $true, $false, $false, $false, $false, $true, 1, $false, "45×29=", 2)
